$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("G3").Value = 65
$ws.Range("F7").Value = 223
$ws.Range("F9").Value = 293
$ws.Range("F11").Value = 857
$ws.Range("F12").Value = 647
$ws.Range("F15").Value = 138
$ws.Range("F18").Value = 2851
$ws.Range("F26").Value = 2426
$ws.Range("F28").Value = 973
$ws.Range("F29").Value = 13
$ws.Range("F31").Value = 266
$ws.Range("F32").Value = 1051
$ws.Range("F35").Value = 277

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 1063
$ws.Range("F5").Value = 1063
$ws.Range("F16").Value = 97
$ws.Range("F18").Value = 973
$ws.Range("F21").Value = 613
$ws.Range("F25").Value = 301
$ws.Range("F27").Value = 3818
$ws.Range("F32").Value = 33
$ws.Range("F34").Value = 142

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2420
$ws.Range("F6").Value = 1004
$ws.Range("F10").Value = 336

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 2420
$ws.Range("G5").Value = 65
$ws.Range("F6").Value = 1004
$ws.Range("F8").Value = 336
$ws.Range("F13").Value = 223
$ws.Range("F16").Value = 293
$ws.Range("F17").Value = 857
$ws.Range("F18").Value = 647
$ws.Range("F19").Value = 1063
$ws.Range("F21").Value = 138
$ws.Range("F24").Value = 2851
$ws.Range("F31").Value = 2426
$ws.Range("F33").Value = 973
$ws.Range("F36").Value = 13
$ws.Range("F37").Value = 97
$ws.Range("F39").Value = 266
$ws.Range("F44").Value = 301
$ws.Range("F45").Value = 301
$ws.Range("F47").Value = 1051
$ws.Range("F49").Value = 142
$ws.Range("F51").Value = 277
